# Update the data table on Sheet1:
#  - Column header "Materi" is replaced by "NIP_NUPTK"
#  - Column "Web"/"Test" content is replaced by numeric NIP/NUPTK-like
#    numbers for rows 2-8, while column D keeps "Test"
#  - The previously unused strings "Materi" and "Web" are dropped since
#    no cell references them any more

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("C1").Value = "NIP_NUPTK"
$ws.Range("D1").Value = "Tema Webinar"

# Data rows 2-8: column C becomes a plain number, column D stays "Test"
$ws.Range("C2").Value = 123
$ws.Range("D2").Value = "Test"

$ws.Range("C3").Value = 234
$ws.Range("D3").Value = "Test"

$ws.Range("C4").Value = 345
$ws.Range("D4").Value = "Test"

$ws.Range("C5").Value = 456
$ws.Range("D5").Value = "Test"

$ws.Range("C6").Value = 567
$ws.Range("D6").Value = "Test"

$ws.Range("C7").Value = 678
$ws.Range("D7").Value = "Test"

$ws.Range("A8").Value = "Habi"
$ws.Range("B8").Value = "IPS"
$ws.Range("C8").Value = 789
$ws.Range("D8").Value = "Test"

# Update the active selection shown in the workbook to C1
$ws.Range("C1").Select()
